# Apply the WBS.xlsx edit: add a "% complete" style column entry and a note
# column entry for the "Nhập hội viên" and "Cập nhật thông tin hội viên" tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11: "Nhập hội viên" task -> mark 100% done, add note
$ws.Range("F11").Value = 1
$ws.Range("F11").NumberFormat = "0%"
$ws.Range("G11").Value = "Có thay đổi giao diện"

# Row 12: "Cập nhật thông tin hội viên" task -> mark 100% done, add note
$ws.Range("F12").Value = 1
$ws.Range("F12").NumberFormat = "0%"
$ws.Range("G12").Value = "Có thay đổi giao diện"

# Match the border formatting used by the rest of row 11 (no top border on G12)
$ws.Range("G12").Borders.Item(8).LineStyle = -4142

# Reset the window scroll position so the sheet view no longer reports a
# scrolled topLeftCell, and move the active selection to C11.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C11").Select()
